$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 47, shifting the existing row 47 (and below) down to row 48.
$ws.Rows.Item(47).Insert()

# Fill the new row 47 with the new data entry.
$ws.Range("A47").Value = 8
$ws.Range("B47").Value = "Terminal La Palmera de La Serena"
$ws.Range("C47").Value = "Coquimbo"
$ws.Range("D47").Value = 44747
$ws.Range("D47").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E47").Value = 4
$ws.Range("F47").Value = 100114007
$ws.Range("G47").Value = "Jengibre"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 440
$ws.Range("K47").Value = 15000
$ws.Range("L47").Value = 16000
$ws.Range("M47").Value = 15500
$ws.Range("N47").Value = "`$/caja 13 kilos"
$ws.Range("O47").Value = "Perú"
$ws.Range("P47").Value = 1192
$ws.Range("Q47").Value = 13
$ws.Range("R47").Value = "Hortaliza"
